# This script reproduces the following content changes against
# Project_1_Update.docx:
#
#  1. Insert a new "Implement rfe feature selection" bullet just before the
#     "Implement sentiment model..." bullet (Work Left list).
#  2. Append " (experimenting and researching online tutorials for good
#     definitions)" to the "Defining a good cutoff for rare words" bullet.
#  3. Append " (feature selection will help with this but I'm concerned
#     about lowered performance" to the "...is currently over 5000" bullet.
#  4. Append " (once again, feature selection)" to the "Slow processing
#     speeds..." bullet.
#  5. Rewrite "Defining a good metric for feature selection" into
#     "Defining a good algorithm for feature selection (currently removing
#     correlated features, will also apply rfe to full dataset" (split
#     across several runs, with "rfe" wrapped in spell-check proofErr
#     markers, same as the rest of the document).
#
# Paragraph indices shift whenever a new paragraph is inserted, so the
# edits below are applied from the bottom of the document upwards -- that
# way every index used is still the index from the *original* document
# layout at the time it is used.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function New-PkgXml([string]$bodyInnerXml) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="' + $wNs + '"><w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Replaces the visible text (everything up to, but excluding, the paragraph
# mark) of paragraph $index with the run/proofErr markup in $innerParaXml,
# while leaving the paragraph's own pPr (style/numbering/rsid/...) intact.
function Set-ParagraphRunXml([int]$index, [string]$innerParaXml) {
    $para = $d.Paragraphs.Item($index)
    $rangeStart = $para.Range.Start
    $rangeEnd = $para.Range.End
    $bodyRange = $d.Range($rangeStart, $rangeEnd - 1)
    $xml = New-PkgXml "<w:p>$innerParaXml</w:p>"
    $bodyRange.InsertXML($xml) | Out-Null
}

# ---------------------------------------------------------------------------
# Change 5 (paragraph 28): "Defining a good metric for feature selection"
# ---------------------------------------------------------------------------
$inner5 = '<w:r><w:t xml:space="preserve">Defining a good </w:t></w:r><w:r><w:t xml:space="preserve">algorithm for feature selection (currently removing correlated features, will also apply </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rfe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to full dataset</w:t></w:r>'
Set-ParagraphRunXml 28 $inner5

# ---------------------------------------------------------------------------
# Change 4 (paragraph 27): "Slow processing speeds..."
# ---------------------------------------------------------------------------
$inner4 = '<w:r><w:t xml:space="preserve">Slow processing speeds due to the high number of features and complexity of the dataset</w:t></w:r><w:r><w:t xml:space="preserve"> (once again, feature selection)</w:t></w:r>'
Set-ParagraphRunXml 27 $inner4

# ---------------------------------------------------------------------------
# Change 3 (paragraph 26): "...is currently over 5000"
# (the unmatched "(" is intentional -- it matches the source commit)
# ---------------------------------------------------------------------------
$inner3 = '<w:r><w:t xml:space="preserve">The number of words, even with removal of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stopwords</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is currently over 5000</w:t></w:r><w:r><w:t xml:space="preserve"> (feature selection will help with this but I&#8217;m concerned about lowered performance</w:t></w:r>'
Set-ParagraphRunXml 26 $inner3

# ---------------------------------------------------------------------------
# Change 2 (paragraph 25): "Defining a good cutoff for rare words"
# ---------------------------------------------------------------------------
$inner2 = '<w:r><w:t xml:space="preserve">Defining a good cutoff for rare words</w:t></w:r><w:r><w:t xml:space="preserve"> (experimenting and researching online tutorials for good definitions)</w:t></w:r>'
Set-ParagraphRunXml 25 $inner2

# ---------------------------------------------------------------------------
# Change 1 (paragraph 13): insert new bullet before "Implement sentiment
# model..."
# ---------------------------------------------------------------------------
$targetPara = $d.Paragraphs.Item(13)
$targetPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item(13)
$inner1 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Implement </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rfe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> feature selection</w:t></w:r>'
$xml1 = New-PkgXml "<w:p>$inner1</w:p>"
$newPara.Range.InsertXML($xml1) | Out-Null

Write-Output "edits applied"
